# Apply the "Updated cryptos list" data refresh (GitHub Actions commit) to the worksheet.
# All Price column (D) values are plain text in the source data (e.g. "30.190.70",
# "0.4720") so we force text formatting before/while assigning them to avoid Excel's
# automatic type coercion turning them into numbers (which would also strip
# significant trailing zeros, e.g. "0.4720" -> 0.472).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Bitcoin
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.190.70'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  -0.47%  '

# Row 3 - Ethereum
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.862.73'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  -0.40%  '

# Row 4 - TetherUSD
$ws.Range('E4').Value = '  -0.18%  '

# Row 5 - BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.98'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +3.28%  '

# Row 6 - USDC
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9996'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -0.16%  '

# Row 7 - XRP
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4720'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  +0.52%  '

# Row 8 - OKB
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.77'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  -0.37%  '

# Row 9 - Cardano
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2858'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  -0.24%  '

# Row 10 - Dogecoin
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06478'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -1.68%  '

# Row 11 - Solana
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.81'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  -4.49%  '

# Row 12 - TRON
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07680'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  -3.87%  '

# Row 13 - WrappedEther
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.861.20'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  -0.64%  '

# Row 14 - Litecoin
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.16'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -2.86%  '

# Row 15 - Polkadot->Polygon (row/row16 swap)
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6825'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  -1.11%  '

# Row 16 - Polygon->Polkadot (row15/row swap)
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.077'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -0.73%  '

# Row 17 - BitcoinCash
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '269.81'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  +0.17%  '

# Row 18 - WrappedBTC
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '30.181.27'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  -0.50%  '

# Row 19 - Avalanche
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.35'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -5.76%  '

# Row 20 - ShibaInu
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007545'
$ws.Range('D20').NumberFormat = 'General'

# Row 21 - Dai
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9994'
$ws.Range('D21').NumberFormat = 'General'

# Row 22 - WrappedliquidstakedEther2.0
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.103.28'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -0.51%  '

# Row 23 - BinanceUSD
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9990'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -0.19%  '

# Row 24 - Uniswap
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.176'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -1.70%  '

# Row 25 - Chainlink
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.101'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -1.89%  '

# Row 26 - Cosmos
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.345'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -0.57%  '

# Row 27 - Monero
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.02'
$ws.Range('D27').NumberFormat = 'General'

# Row 28 - EthereumClassic
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.73'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -0.85%  '

# Row 29 - LidoDAOToken
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.885'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -3.37%  '

# Row 30 - Toncoin
$ws.Range('E30').Value = '  +0.62%  '

# Row 31 - Stellar
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09856'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -0.11%  '

# Row 32 - PancakeSwap
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.507'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +3.28%  '

# Row 33 - Filecoin
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.234'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -2.76%  '

# Row 34 - InternetComputer(DFINITY)
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.998'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  -1.52%  '

# Row 35 - Hedera
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.04708'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -0.17%  '

# Row 36 - ARBITRUM
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.111'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  -2.19%  '

# Row 37 - ImmutableX
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6864'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -2.26%  '

# Row 38 - HuobiToken
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.708'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -1.25%  '

# Row 39 - VeChain
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01830'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -2.75%  '

# Row 40 - MXToken
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.725'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -3.73%  '

# Row 41 - FraxShare
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.387'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  +2.15%  '

# Row 42 - Aave
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '70.25'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -2.60%  '

# Row 43 - PaxDollar
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9992'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -0.16%  '

# Row 44 - TrustWalletToken
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8351'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -0.87%  '

# Row 45 - RenderToken
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.891'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -3.52%  '

# Row 46 - Quant
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.14'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -0.94%  '

# Row 47 - TheSandbox
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4072'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  -2.58%  '

# Row 48 - EnergySwap
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.248'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +1.45%  '

# Row 49 - Maker
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '927.11'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  +0.96%  '

# Row 50 - Aptos
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.943'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -2.10%  '

# Row 51 - Elrond
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '34.42'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -0.37%  '
